{"js": "// Readme.docx: \"Building the sample\" section, Project Scarlett paragraph.\n// Before: \"If using Project Scarlett, set the active solution platform to Gaming.Xbox.Scarlett.x64.\"\n// After:  \"If using an Xbox Series X|S devkit, set the active solution platform to Gaming.Xbox.Scarlett.x64.\"\nconst oldPhrase = \"Project Scarlett\";\nconst newPhrase = \"an Xbox Series X|S devkit\";\n\nconst results = context.document.body.search(oldPhrase, { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(`Could not find the text \"${oldPhrase}\" to replace.`);\n}\n\n// Replace just the first (and expected only) occurrence, preserving the\n// surrounding run formatting (the match carries no special formatting here).\nresults.items[0].insertText(newPhrase, \"Replace\");\nawait context.sync();\n", "ps1": "# Readme.docx: \"Building the sample\" section, Project Scarlett paragraph.\n# Before: \"If using Project Scarlett, set the active solution platform to Gaming.Xbox.Scarlett.x64.\"\n# After:  \"If using an Xbox Series X|S devkit, set the active solution platform to Gaming.Xbox.Scarlett.x64.\"\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Project Scarlett\"\n$find.Replacement.Text = \"an Xbox Series X|S devkit\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.Forward = $true\n$find.Wrap = 1  # wdFindContinue\n\n# wdReplaceAll = 2\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n"}
